$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.436.80"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.05"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9973"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9980"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.55"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3757"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.57"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3579"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08186"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.232"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9983"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.34"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.505"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.340"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001224"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.619.52"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.97"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06941"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.703"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.45"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9976"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.48"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.423.17"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.520"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.117"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.14"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.34"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.181"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.91"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.802.40"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.099"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.648"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.62"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.028"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -10.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02749"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08763"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.964"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06961"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.52"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6983"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.330"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.64"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6456"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9978"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.299"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.961"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07923"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.57"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.178"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.85%  "
